# The second paragraph of the document contains the M2Doc field text
# "{m:null.rotate(null)}" spread across several runs. Two of those runs
# need to be split in two (without altering the visible text or
# formatting of the paragraph):
#   "{m"  -> "{" + "m"
#   ")}"  -> ")" + "}"
#
# Re-assigning a Range's own FormattedText back onto itself is enough to
# force Word to break the underlying run at the range boundaries while
# keeping the original run formatting (including leaving runs with no
# explicit rPr untouched), which is exactly what the target edit needs.

$d = $word.ActiveDocument
$p = $d.Paragraphs(2)

$paraStart = $p.Range.Start
$paraEnd = $p.Range.End

# --- Split "{m" into "{" and "m" -------------------------------------
$openBrace = $d.Range($paraStart, $paraStart + 1)
$openBrace.FormattedText = $openBrace.FormattedText

# --- Split ")}" into ")" and "}" --------------------------------------
# paraEnd points just past the paragraph mark, so the closing "}" sits
# two characters before it, and ")" one more character before that.
$closeParen = $d.Range($paraEnd - 3, $paraEnd - 2)
$closeParen.FormattedText = $closeParen.FormattedText
